$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, as scraped from the updated cryptos feed.
$updates = [ordered]@{
    'D2' = '45.339.86'
    'E2' = '  +6.46%  '
    'D3' = '2.382.41'
    'E3' = '  +3.65%  '
    'E4' = '  -0.28%  '
    'D5' = '111.86'
    'E5' = '  +6.88%  '
    'D6' = '317.26'
    'E6' = '  +1.31%  '
    'D7' = '0.639'
    'E7' = '  +2.35%  '
    'E8' = '  +0.03%  '
    'D9' = '0.632'
    'E9' = '  +4.64%  '
    'D10' = '42.34'
    'E10' = '  +7.79%  '
    'D11' = '0.0930'
    'E11' = '  +2.73%  '
    'D12' = '8.71'
    'E12' = '  +4.90%  '
    'D13' = '1.02'
    'E13' = '  +2.42%  '
    'D14' = '0.108'
    'E14' = '  +0.48%  '
    'D15' = '15.83'
    'E15' = '  +4.18%  '
    'D16' = '2.742.05'
    'E16' = '  +3.56%  '
    'D17' = '2.375.09'
    'E17' = '  +3.28%  '
    'D18' = '45.296.85'
    'E18' = '  +5.92%  '
    'D19' = '7.65'
    'E19' = '  +4.37%  '
    'E20' = '  +3.02%  '
    'D21' = '13.15'
    'E21' = '  -3.40%  '
    'D22' = '75.25'
    'E22' = '  +2.44%  '
    'E23' = '  +1.87%  '
    'D24' = '269.77'
    'E24' = '  +1.59%  '
    'D25' = '2.35'
    'E25' = '  +7.25%  '
    'E26' = '  -0.44%  '
    'D27' = '7.69'
    'E27' = '  +7.59%  '
    'D28' = '11.30'
    'E28' = '  +4.89%  '
    'E29' = '  -0.16%  '
    'D30' = '39.65'
    'E30' = '  +9.69%  '
    'D31' = '22.92'
    'D32' = '0.0934'
    'E32' = '  +7.88%  '
    'D33' = '169.36'
    'E33' = '  +2.46%  '
    'D34' = '2.98'
    'E34' = '  +16.30%  '
    'E35' = '  +1.98%  '
    'E36' = '  +3.40%  '
    'D37' = '4.83'
    'E37' = '  +6.71%  '
    'D38' = '0.0367'
    'E38' = '  +4.78%  '
    'D39' = '3.03'
    'E39' = '  +9.78%  '
    'D40' = '3.92'
    'E40' = '  +4.69%  '
    'E41' = '  +8.45%  '
    'D42' = '106.94'
    'E42' = '  +5.52%  '
    'D43' = '13.87'
    'E43' = '  +13.86%  '
    'E44' = '  +6.15%  '
    'D45' = '71.87'
    'E45' = '  +3.46%  '
    'E46' = '  -0.26%  '
    'D47' = '119.15'
    'E47' = '  +6.54%  '
    'B48' = 'THORChain'
    'C48' = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
    'D48' = '5.70'
    'E48' = '  +9.39%  '
    'B49' = 'MinaProtocolToken'
    'C49' = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
    'D49' = '1.65'
    'E49' = '  +19.58%  '
    'D50' = '79.60'
    'E50' = '  -0.63%  '
    'D51' = '0.221'
    'E51' = '  +16.62%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Every updated cell holds plain text (inline strings), including column D
    # prices such as "45.339.86" that look numeric. Force text formatting first
    # so Excel does not silently coerce the assignment into a float and drop
    # the grouping dots / trailing zeros.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
